$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- Task "CU Consultar proximos pagos de profesores." (row 8) is now Done ---
# Status column F8: "En proceso" -> "Hecho"
$ws.Range("F8").Value = "Hecho"

# Register 3 consumed hours on day 3 (column N), which cascades the
# "Rest." shared formulas across the remaining days and flips the
# totals (AZ8 / BA8) accordingly.
$ws.Range("N8").Value = 3

# --- Move the current selection to Z13 (bottomRight pane) ---
$ws.Range("Z13").Select()

# --- Re-merge the day-total header cells; this reshuffles the stored
# mergeCells order the same way Excel does when the merged range is
# touched (unmerge, then re-merge each piece). ---
$ws.Range("H4:BA4").UnMerge()
$ws.Range("AZ4:BA4").Merge()
$ws.Range("AO4:AP4").Merge()
$ws.Range("AR4:AS4").Merge()
$ws.Range("AU4:AV4").Merge()
$ws.Range("AX4:AY4").Merge()
$ws.Range("AL4:AM4").Merge()
$ws.Range("H4:I4").Merge()
$ws.Range("K4:L4").Merge()
$ws.Range("N4:O4").Merge()
$ws.Range("Q4:R4").Merge()
$ws.Range("T4:U4").Merge()
$ws.Range("W4:X4").Merge()
$ws.Range("Z4:AA4").Merge()
$ws.Range("AC4:AD4").Merge()
$ws.Range("AF4:AG4").Merge()
$ws.Range("AI4:AJ4").Merge()

$wb.Save()
